# Fixed typo in sample files: "LicenseTier*" -> "LicenceTier*"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WMT_Extract")

$ws.Range("W1").Value = "LicenceTier0"
$ws.Range("X1").Value = "LicenceTierD2"
$ws.Range("Y1").Value = "LicenceTierD1"
$ws.Range("Z1").Value = "LicenceTierC2"
$ws.Range("AA1").Value = "LicenceTierC1"
$ws.Range("AB1").Value = "LicenceTierB2"
$ws.Range("AC1").Value = "LicenceTierB1"
$ws.Range("AD1").Value = "LicenceTierA"

$ws.Range("AD2").Select()
